# "Dados apresentação 22-08.xlsx" / Planilha1
# Commit: "Updated test suite results ... MIRO transaction" edit —
# the MIRO-linked purchase order / invoice numbers (cols P = "NV CONTRATO",
# Q = "NV PEDIDO") were bumped to the next set of numbers for rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 4600244316
$ws.Range("Q2").Value = 4503342084
$ws.Range("P3").Value = 4600244317
$ws.Range("Q3").Value = 4503342085

# Match the author's final selection/scroll state (selection collapsed to Q3,
# view scrolled right so column I is the leftmost visible column).
$ws.Range("Q3").Select()
$excel.ActiveWindow.ScrollColumn = 9
